# Update the "Förändrad" (Changed) date column (C) for rows 2-10:
# change the stored serial date value from 45171 to 45172 (i.e. +1 day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
